$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:E (leg trial 1,2) get overwritten with the values that were
# previously in columns O, R, AN, AQ (leg trial 16,20) - i.e. the old
# "1,2,3,4" trial columns are replaced by duplicates of the 16/20 trial
# columns, effectively deleting the original 1/2/3/4 values.

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 105.98916852820224
$ws.Range("C2").Value = 104.74085481089448
$ws.Range("D2").Value = 107.76606483851549
$ws.Range("E2").Value = 107.10477109939815

$ws.Range("B3").Value = 104.88524901633632
$ws.Range("C3").Value = 107.8360232974745
$ws.Range("D3").Value = 108.64319819792583
$ws.Range("E3").Value = 106.41734465713107

# Update the selected range to reflect the narrower region of interest.
$ws.Range("B1:E3").Select()
